$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 243 (C, D only - A, B unchanged)
$ws.Range("C243").Value = 31.28546357154846
$ws.Range("D243").Value = 27.65625

# Update row 244 (A, B, C, D all change)
$ws.Range("A244").Value = 82.45
$ws.Range("B244").Value = 0.16
$ws.Range("C244").Value = 41.23247504234314
$ws.Range("D244").Value = 38.109375

# Update row 245 (A, B, C, D all change)
$ws.Range("A245").Value = 39.78
$ws.Range("B245").Value = 0.16
$ws.Range("C245").Value = 37.43526697158813
$ws.Range("D245").Value = 36.25

# Update row 246 (C, D only - A, B unchanged)
$ws.Range("C246").Value = 34.47661805152893
$ws.Range("D246").Value = 33.0625

# Update row 247 (C, D only - A, B unchanged)
$ws.Range("C247").Value = 27.1760106086731
$ws.Range("D247").Value = 25.6875

# Remove rows 248-253 (speed run data trimmed)
$ws.Rows("248:253").Delete()
